$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (soln) ---
# Update row 2 values
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = 2
$ws1.Range("C2").Value = 0
$ws1.Range("D2").Value = 77
$ws1.Range("E2").Value = 300

# Unmerge cells A2:A4 before deleting rows
$ws1.Range("A2:A4").UnMerge()

# Delete rows 3 and 4 entirely
$ws1.Rows("3:4").Delete()

# --- Sheet2 (investmentsoln) ---
$ws2.Range("A2").Value = 2
$ws2.Range("A3").Value = 3
$ws2.Range("A4").Value = 4
$ws2.Range("B4").Value = 0

# Copy the style of A3 to A4 so it matches the bordered/bold style
$ws2.Range("A3").Copy()
$ws2.Range("A4").PasteSpecial(-4122)
$ws2.Range("A4").Value = 4
